$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.455.31'
$ws.Range('E2').Value = '  -1.18%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.734.93'
$ws.Range('E3').Value = '  -1.32%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.47%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '322.49'
$ws.Range('E5').Value = '  +0.23%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +0.46%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4609'
$ws.Range('E7').Value = '  +8.98%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3518'
$ws.Range('E8').Value = '  -2.91%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '41.89'
$ws.Range('E9').Value = '  -1.25%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07349'
$ws.Range('E10').Value = '  -1.56%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.081'
$ws.Range('E11').Value = '  -0.23%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.005'
$ws.Range('E12').Value = '  +0.48%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.43'
$ws.Range('E13').Value = '  -1.24%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.908'
$ws.Range('E14').Value = '  -2.53%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.041'
$ws.Range('E15').Value = '  -3.40%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.739.08'
$ws.Range('E16').Value = '  -0.66%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '90.99'
$ws.Range('E17').Value = '  +0.07%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001051'
$ws.Range('E18').Value = '  -0.23%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06393'
$ws.Range('E19').Value = '  +0.49%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.45%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.63'
$ws.Range('E21').Value = '  -2.09%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.717'
$ws.Range('E22').Value = '  -3.54%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '27.524.51'
$ws.Range('E23').Value = '  -0.97%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.08'
$ws.Range('E24').Value = '  -0.99%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.107'
$ws.Range('E25').Value = '  +0.54%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.69'
$ws.Range('E26').Value = '  +3.44%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.85'
$ws.Range('E27').Value = '  -1.76%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.941.03'
$ws.Range('E28').Value = '  -0.59%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '124.52'
$ws.Range('E29').Value = '  +0.53%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.034'
$ws.Range('E30').Value = '  -4.42%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.042'
$ws.Range('E31').Value = '  -6.19%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09167'
$ws.Range('E32').Value = '  +3.58%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.666'
$ws.Range('E33').Value = '  -0.47%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.414'
$ws.Range('E34').Value = '  -2.38%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02264'
$ws.Range('E35').Value = '  -1.05%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '11.59'
$ws.Range('E36').Value = '  -4.96%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05985'
$ws.Range('E37').Value = '  -0.90%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2062'
$ws.Range('E38').Value = '  -1.51%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.909'
$ws.Range('E39').Value = '  -0.58%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6235'
$ws.Range('E40').Value = '  -1.18%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.183'
$ws.Range('E41').Value = '  +0.55%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.376'
$ws.Range('E42').Value = '  -1.55%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.693'
$ws.Range('E43').Value = '  -2.23%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.00'
$ws.Range('E44').Value = '  -1.75%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.698'
$ws.Range('E45').Value = '  +0.41%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5796'
$ws.Range('E46').Value = '  -1.01%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '121.58'
$ws.Range('E47').Value = '  -1.08%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.917'
$ws.Range('E48').Value = '  -3.08%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06830'
$ws.Range('E49').Value = '  +0.26%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.120'
$ws.Range('E50').Value = '  -4.81%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '71.13'
$ws.Range('E51').Value = '  -3.37%  '
